$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.953.29"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "'1.549.26"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'205.85"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'0.484"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "'21.45"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "'1.769.56"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "'1.548.44"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "'3.70"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "'0.513"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "'26.928.44"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "'61.58"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "'214.26"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "'7.22"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "'4.03"
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("D23").Value = "'9.17"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").Value = "'6.65"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "'14.83"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'0.0460"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").Value = "'1.371.47"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").Value = "'0.967"
$ws.Range("E36").Value = "  +4.70%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "'0.518"
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("D40").Value = "'0.806"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "'0.986"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("D45").Value = "'63.57"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("B47").Value = "mCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D47").Value = "'2.26"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "'1.683.43"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "'86.19"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").Value = "'0.0507"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").Value = "'0.0952"
$ws.Range("E51").Value = "  -0.14%  "
